# Updating filtered feeds from workflow
# Appends two new rows (94, 95) to the "Filtered Feeds" sheet, one for each
# source domain (genomeweb.com / 360dx.com) of a new article about the
# Merck/Agilent PD-L1 ovarian cancer companion diagnostic, mirroring the
# existing link/keywords/title layout used by every prior row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$title = "Merck Immunotherapy, Agilent CDx Net FDA Approval for PD-L1-Positive Ovarian Cancer"
$keywords = "CDx"

$link94 = "https://www.genomeweb.com/cancer/merck-immunotherapy-agilent-cdx-net-fda-approval-pd-l1-positive-ovarian-cancer"
$link95 = "https://www.360dx.com/cancer/merck-immunotherapy-agilent-cdx-net-fda-approval-pd-l1-positive-ovarian-cancer"

# Row 94
$ws.Range("A94").Value2 = $link94
$ws.Range("B94").Value2 = $keywords
$ws.Range("C94").Value2 = $title
$ws.Hyperlinks.Add($ws.Range("A94"), $link94) | Out-Null
$ws.Range("A94").Style = "Hyperlink"

# Row 95
$ws.Range("A95").Value2 = $link95
$ws.Range("B95").Value2 = $keywords
$ws.Range("C95").Value2 = $title
$ws.Hyperlinks.Add($ws.Range("A95"), $link95) | Out-Null
$ws.Range("A95").Style = "Hyperlink"
